# Apply the "Dep Ed Closures" content refresh (Thu 27 Aug -> Fri 28 Aug snapshot).
# The sheet is a single column (A) list of scraped page fragments; the refreshed
# snapshot reorders/edits several list items which shifts many row numbers down
# column A. We set each affected cell explicitly to its final value below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = 'On this pageCurrent school and early childhood service, TAFE closures and relocations:Bus service cancellations or alterationsCurrent school and early childhood service, TAFE closures and relocations for Friday 28 August, (as at 10:45am, 28August)South-Eastern Victoria RegionEarly childhood services'
$ws.Cells.Item(38, 1).Value = 'li: Cranbourne Community House CRANBOURNE'
$ws.Cells.Item(39, 1).Value = 'li: David Collings Centre MORNINGTON'
$ws.Cells.Item(40, 1).Value = 'li: Fernwood Fitness Clayton CLAYTON'
$ws.Cells.Item(41, 1).Value = 'li: Fernwood Fitness Narre Warren NARRE WARREN'
$ws.Cells.Item(42, 1).Value = 'li: GEKA Bentleigh Kindergarten BENTLEIGH'
$ws.Cells.Item(43, 1).Value = 'li: GEKA Caulfield South Kindergarten CAULFIELD SOUTH'
$ws.Cells.Item(44, 1).Value = 'li: GEKA Centre Road Kindergarten BENTLEIGH EAST'
$ws.Cells.Item(45, 1).Value = 'li: GEKA Glover Street Kindergarten BENTLEIGH EAST'
$ws.Cells.Item(46, 1).Value = 'li: GEKA McKinnon Kindergarten MCKINNON'
$ws.Cells.Item(47, 1).Value = 'li: GEKA Murrumbeena Kindergarten MURRUMBEENA'
$ws.Cells.Item(48, 1).Value = 'li: GEKA Orrong Road Kindergarten ELSTERNWICK'
$ws.Cells.Item(49, 1).Value = 'li: Genesis Traralgon TRARALGON'
$ws.Cells.Item(50, 1).Value = 'li: GESAC Childcare BENTLEIGH EAST'
$ws.Cells.Item(51, 1).Value = 'li: Glen Iris Creche GLEN IRIS'
$ws.Cells.Item(52, 1).Value = 'li: Goodlife Chelsea Heights CHELSEA HEIGHTS'
$ws.Cells.Item(53, 1).Value = 'li: Goodlife Fountain Gate NARRE WARREN'
$ws.Cells.Item(54, 1).Value = 'li: Goodlife Karingal KARINGAL'
$ws.Cells.Item(55, 1).Value = 'li: Hallam Community Centre Inc HALLAM'
$ws.Cells.Item(56, 1).Value = 'li: Hoa Nghiem Primary School SPRINGVALE SOUTH'
$ws.Cells.Item(57, 1).Value = 'li: Lakes Aquadome Creche LAKES ENTRANCE'
$ws.Cells.Item(58, 1).Value = 'li: Little Beacons Learning Centre - Berwick BERWICK'
$ws.Cells.Item(59, 1).Value = 'li: Little Pelicans - Pelican Park Recreation Centre HASTINGS'
$ws.Cells.Item(60, 1).Value = 'li: Merkaz Bentleigh BENTLEIGH'
$ws.Cells.Item(61, 1).Value = 'li: Moe Neighbourhood House MOE'
$ws.Cells.Item(62, 1).Value = 'li: Paisley Park Early Learning Centre Cranbourne CRANBOURNE'
$ws.Cells.Item(63, 1).Value = 'li: Pavillion - Frankston & District Netball Association Inc FRANKSTON'
$ws.Cells.Item(64, 1).Value = 'li: Peninsula Aquatic Recreation Centre FRANKSTON'
$ws.Cells.Item(65, 1).Value = 'li: Play Zone - Highett HIGHETT'
$ws.Cells.Item(66, 1).Value = 'li: Playtime Warragul WARRAGUL'
$ws.Cells.Item(67, 1).Value = 'li: Sacre Coeur Oshclub GLEN IRIS'
$ws.Cells.Item(68, 1).Value = 'li: Sandybeach Centre SANDRINGHAM'
$ws.Cells.Item(69, 1).Value = 'li: Singleton Preschool ENDEAVOUR HILLS'
$ws.Cells.Item(70, 1).Value = 'li: Smaland Springvale Ikea SPRINGVALE'
$ws.Cells.Item(71, 1).Value = 'li: Somerville Recreation Centre Child Care SOMERVILLE'
$ws.Cells.Item(72, 1).Value = 'li: Toorak Primary School OSHClub TOORAK'
$ws.Cells.Item(73, 1).Value = 'li: Upper Beaconsfield Community Early Learning Centre BEACONSFIELD UPPER'
$ws.Cells.Item(74, 1).Value = 'Schools closedThe Department has been advised of the following school closures:'
$ws.Cells.Item(75, 1).Value = 'li: Korowa Anglican Girl''s School, GLEN IRIS'
$ws.Cells.Item(76, 1).Value = 'li: Lighthouse Christian College, CRANBOURNE'
$ws.Cells.Item(77, 1).Value = 'TAFETheDepartment hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood services'
$ws.Cells.Item(78, 1).Value = 'li: The Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(79, 1).Value = 'li: Aqualink Box Hill Creche BOX HILL'
$ws.Cells.Item(80, 1).Value = 'li: Aqualink Nunawading Creche FOREST HILL'
$ws.Cells.Item(81, 1).Value = 'li: Aquarena Childcare DONCASTER'
$ws.Cells.Item(82, 1).Value = 'li: Billanook College - Early Learning Program MOOROOLBARK'
$ws.Cells.Item(83, 1).Value = 'li: Camp Australia - Balwyn North Primary School OSHC BALWYN NORTH'
$ws.Cells.Item(84, 1).Value = 'li: Camp Australia - Boronia K-12 College OSHC BORONIA'
$ws.Cells.Item(85, 1).Value = 'li: Camp Australia - Camberwell Boys Grammar Junior School OSHC CANTERBURY'
$ws.Cells.Item(86, 1).Value = 'li: Camp Australia - Mountain Gate Primary School OSHC FERNTREE GULLY'
$ws.Cells.Item(87, 1).Value = 'li: Camp Australia - St Dominics Primary School OSHC CAMBERWELL EAST'
$ws.Cells.Item(88, 1).Value = 'li: Camp Australia - St Jude the Apostle School OSHC SCORESBY'
$ws.Cells.Item(89, 1).Value = 'li: Camp Australia - St Timothy''s School Vermont OSHC VERMONT'
$ws.Cells.Item(90, 1).Value = 'li: Camp Australia - Strathcona Baptist Girls Junior School OSHC CANTERBURY'
$ws.Cells.Item(91, 1).Value = 'li: Camp Australia - Trinity Grammar OSHC KEW'
$ws.Cells.Item(92, 1).Value = 'li: Camp Australia - Yarra Valley Grammar School OSHC RINGWOOD'
$ws.Cells.Item(93, 1).Value = 'li: Carey OSHClub Kew KEW'
$ws.Cells.Item(94, 1).Value = 'li: Cathedral College Wangaratta Outside School Hours Care WANGARATTA'
$ws.Cells.Item(95, 1).Value = 'li: Croydon Hills Primary School OSHC Program CROYDON HILLS'
$ws.Cells.Item(96, 1).Value = 'li: Ferntree Gully 3 Year Old Preschool FERNTREE GULLY'
$ws.Cells.Item(97, 1).Value = 'li: Fitness First Doncaster (Playzone) DONCASTER'
$ws.Cells.Item(98, 1).Value = 'li: Flamingo Community Group WANTIRNA SOUTH'
$ws.Cells.Item(99, 1).Value = 'li: Genesis Creche Wantirna BORONIA'
$ws.Cells.Item(100, 1).Value = 'li: Hawthorn Early Years HAWTHORN'
$ws.Cells.Item(101, 1).Value = 'li: Highmount Kindergarten MOUNT WAVERLEY'
$ws.Cells.Item(102, 1).Value = 'li: Insight Early Learning - Burwood Brickworks BURWOOD EAST'
$ws.Cells.Item(103, 1).Value = 'li: JOYFUL LEARNING EARLY LEARNING CENTRE PTY LTD MITCHAM'
$ws.Cells.Item(104, 1).Value = 'li: Kew Neighbourhood Learning Centre KEW'
$ws.Cells.Item(206, 1).Value = 'li: Pentland Afterschool Care Program DARLEY'
$ws.Cells.Item(207, 1).Value = 'li: Shine Early Learning St Albans ST ALBANS'
$ws.Cells.Item(208, 1).Value = 'li: Shuter Street Occasional Care MOONEE PONDS'
$ws.Cells.Item(209, 1).Value = 'li: South Kingsville Community Centre SOUTH KINGSVILLE'
$ws.Cells.Item(210, 1).Value = 'li: Springside Occasional Care CAROLINE SPRINGS'
$ws.Cells.Item(211, 1).Value = 'li: Sunshine Leisure Centre SUNSHINE'
$ws.Cells.Item(212, 1).Value = 'li: Warracknabeal Neighbourhood House Take a Break Child Care Centre WARRACKNABEAL'
$ws.Cells.Item(213, 1).Value = 'li: Willaura Primary School OSHC ARARAT'
$ws.Cells.Item(214, 1).Value = 'li: Winchelsea Primary School Winchelsea'
$ws.Cells.Item(215, 1).Value = 'li: Yarraville Community Centre YARRAVILLE'
$ws.Cells.Item(216, 1).Value = 'Schools closed'
$ws.Cells.Item(217, 1).Value = 'li: The Department has been advised of the following school closures:'
$ws.Cells.Item(218, 1).Value = 'Al Taqwa College, TRUGANINAParkville College (Malmsbury campus), PARKVILLEParkville College (Collingwood campus), COLLINGWOODWarringa Park School (Warringa Crescent campus), HOPPERS CROSSINGTAFE'
$ws.Cells.Item(219, 1).Value = 'li: The Department has not been advised of any TAFE closures.'
$ws.Cells.Item(220, 1).Value = 'Bus service cancellations or alterations'
$ws.Cells.Item(221, 1).Value = 'li: For Term 3 2020, schools bus services will continue to be provided to support student travel to schools where needed.'
$ws.Cells.Item(222, 1).Value = 'li: With learning from home arrangements, local principals are authorised to cancel or suspend school buses where not required.'
$ws.Cells.Item(223, 1).Value = 'li: Families and students are encouraged to contact their school directly to enquire if their bus service is continuing to be provided.'
$ws.Cells.Item(224, 1).Value = 'li: Find more about bus services:'
$ws.Cells.Item(225, 1).Value = 'School bus services in Term 3 (for schools)School bus services in Term 3 (for parents)VicRoadsRegional roadsBack to top'
$ws.Cells.Item(226, 1).Value = 'li: Last Update: 28 August 2020'
$ws.Cells.Item(227, 1).Value = 'In this section:- Example left hand nav using current codeHealth, wellbeing and safetyWebsite navigationFor parentsEarly childhoodSchoolsTAFE and trainingThe Department<li>                                 <a href="/about/educationstate">Education State                                                                  </a>                             </li>Support linksHelp in your languageAccessibilityPrivacyCopyrightDisclaimerContact linksContact usReport a website issue'
$ws.Cells.Item(228, 1).Value = 'li: State Government of Victoria, Australia © 2019'
$ws.Cells.Item(229, 1).Value = 'li: We respectfully acknowledge the Traditional Owners of country throughout Victoria and pay respect to the ongoing living cultures of First Peoples.'
$ws.Cells.Item(230, 1).Value = 'li: Our website uses a free tool to translate into other languages. This tool is a guide and may not be accurate. For more, see:'
$ws.Cells.Item(231, 1).Value = 'Information in your languageOld Search Code PlaceholdCludo Code for www.education.vic.gov.auHotjar Tracking Code for www.education.vic.gov.auSR-1181393'
